$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 89.5
$ws.Range("B2").Value = 3.721558813185679
$ws.Range("C2").Value = 91.3
$ws.Range("D2").Value = 2.934280150224242
$ws.Range("E2").Value = 89.8
$ws.Range("F2").Value = 2.638181191654584
$ws.Range("G2").Value = 91.40000000000001
$ws.Range("H2").Value = 3.006659275674582
$ws.Range("I2").Value = 90.59999999999999
$ws.Range("J2").Value = 3.611094017053557
$ws.Range("K2").Value = 90.90000000000001
$ws.Range("L2").Value = 3.389690251335658
$ws.Range("M2").Value = 94.09999999999999
$ws.Range("N2").Value = 2.808914381037628
$ws.Range("O2").Value = 94.09999999999999
$ws.Range("P2").Value = 1.3
$ws.Range("Q2").Value = 92.7
$ws.Range("R2").Value = 2.865309756378881
$ws.Range("S2").Value = 95.40000000000001
$ws.Range("T2").Value = 2.2
$ws.Range("U2").Value = 95.40000000000001
$ws.Range("V2").Value = 3.072458299147443
$ws.Range("W2").Value = 95
$ws.Range("X2").Value = 4
$ws.Range("Y2").Value = 95.2
$ws.Range("Z2").Value = 3.218695387886216
$ws.Range("AA2").Value = 94.40000000000001
$ws.Range("AB2").Value = 4.223742416388575
$ws.Range("AC2").Value = 95
$ws.Range("AD2").Value = 2.607680962081059
$ws.Range("AE2").Value = 97.09999999999999
$ws.Range("AF2").Value = 5.467174773134658
$ws.Range("AG2").Value = 95.90000000000001
$ws.Range("AH2").Value = 3.726929030716844
$ws.Range("AI2").Value = 97
$ws.Range("AJ2").Value = 2.863564212655271
$ws.Range("AK2").Value = 97.40000000000001
$ws.Range("AL2").Value = 2.244994432064364
$ws.Range("AM2").Value = 97.8
$ws.Range("AN2").Value = 4.1182520563948
$ws.Range("AO2").Value = 98.8
$ws.Range("AP2").Value = 2.675817632051931
$ws.Range("AQ2").Value = 98
$ws.Range("AR2").Value = 3.16227766016838
$ws.Range("AS2").Value = 95.7
$ws.Range("AT2").Value = 2.491987158875422
$ws.Range("AU2").Value = 96.8
$ws.Range("AV2").Value = 3.059411708155671
$ws.Range("AW2").Value = 95.7
$ws.Range("AX2").Value = 4.670117771534247
$ws.Range("AY2").Value = 94.5
$ws.Range("AZ2").Value = 3.640054944640259
$ws.Range("BA2").Value = 94.90000000000001
$ws.Range("BB2").Value = 4.253234063627347
$ws.Range("BC2").Value = 93.2
$ws.Range("BD2").Value = 3.841874542459709
$ws.Range("BE2").Value = 93.3
$ws.Range("BF2").Value = 5.080354318352215
$ws.Range("BG2").Value = 92
$ws.Range("BH2").Value = 4.604345773288535
$ws.Range("BI2").Value = 2832.9